$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.9494948151808758, 0.7473989913805523, 1.1031716645083636, 3.2278711134822036)"
$ws.Range("C2").Value = "JSU(-1.421257629543053, 1.212392539523929, 2.323854466628398, 4.605743585923996)"
$ws.Range("D2").Value = "JSU(-0.9764156446533655, 1.0993949661274143, 0.5082753994897704, 2.488298352224855)"
$ws.Range("E2").Value = "NCT(3.044666318331854, 1.6387985353753538, -0.016454828117430403, 4.7913157623756)"
